# Mise à jour de l'application
# Adds a new attendance column (training day dated 2025-10-29, Excel serial
# 45959) right after the last existing day column (BR, serial 45958) on the
# single worksheet "Feuil1".
#
# Column BS1 gets the new date; BS2:BS29 (except row 12, whose player left
# the roster earlier and whose data already stops at column AX) get the
# attendance code for that day, copying the formatting of the prior day's
# column (BR) so cell styles line up with the rest of the sheet. The summary
# formulas in columns B:J are simple COUNTA/COUNTIF ranges that already
# stretch far beyond column BS, so they recalculate automatically once the
# new cells exist.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New day header (column BS, row 1) ------------------------------------
# (Value is set before the format-only paste so the workbook's dependency
# tracker reliably marks the summary formulas dirty for recalculation.)
$ws.Range("BS1").Value = 45959
$ws.Range("BR1").Copy()
$ws.Range("BS1").PasteSpecial(-4122)  # xlPasteFormats

# --- Attendance codes for the new day, row by row --------------------------
$attendance = @{
    2  = "P"
    3  = "P"
    4  = "P"
    5  = "P"
    6  = "B"
    7  = "P"
    8  = "B"
    9  = "P"
    10 = "P"
    11 = "P"
    13 = "B"
    14 = "P"
    15 = "P"
    16 = "P"
    17 = "P"
    18 = "P"
    19 = "P"
    20 = "P"
    21 = "B"
    22 = "P"
    23 = "RH"
    24 = "P"
    25 = "P"
    26 = "P"
    27 = "P"
    28 = "P"
    29 = "P"
}

foreach ($row in $attendance.Keys) {
    $srcCell = $ws.Cells.Item($row, 70)  # BR<row>
    $dstCell = $ws.Cells.Item($row, 71)  # BS<row>
    $dstCell.Value = $attendance[$row]
    $srcCell.Copy()
    $dstCell.PasteSpecial(-4122)  # xlPasteFormats
}

$excel.CutCopyMode = $false

# Row 12's player roster entry ends at column AX, well before BR/BS, so it
# intentionally receives no new cell (matches the source data).

# --- Refresh view state: keep the frozen first column, move the selection
# to the new edge of the data like the original edit session ended up. -----
$ws.Range("BV17").Select()

# Force the COUNTA/COUNTIF summary formulas (columns B:J) to recompute now
# that the new day's data exists.
$excel.Calculate()
